$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 43193

$ws.Range("D8:G8").ClearContents()
$ws.Range("A8:C8").Style = "Bueno"

$ws.Range("G18").Select()
